$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.664.99"
$ws.Range("E2").Value = "  +1.28%  "
$ws.Range("D3").Value = "3.400.72"
$ws.Range("E3").Value = "  +0.37%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'561.29"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("D6").Value = "'176.27"
$ws.Range("E6").Value = "  +0.65%  "
$ws.Range("D7").Value = "'0.633"
$ws.Range("E7").Value = "  +0.84%  "
$ws.Range("D8").Value = "3.393.04"
$ws.Range("E8").Value = "  +0.46%  "
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("E10").Value = "  +4.99%  "
$ws.Range("D11").Value = "'0.637"
$ws.Range("E11").Value = "  +0.86%  "
$ws.Range("D12").Value = "'53.54"
$ws.Range("E12").Value = "  -1.80%  "
$ws.Range("D13").Value = "'0.0000278"
$ws.Range("E13").Value = "  +0.47%  "
$ws.Range("E14").Value = "  +0.93%  "
$ws.Range("D15").Value = "3.935.73"
$ws.Range("E15").Value = "  -0.02%  "
$ws.Range("D16").Value = "'18.31"
$ws.Range("E16").Value = "  +0.24%  "
$ws.Range("D17").Value = "3.405.81"
$ws.Range("E17").Value = "  +0.57%  "
$ws.Range("E18").Value = "  +1.08%  "
$ws.Range("D19").Value = "65.587.73"
$ws.Range("E19").Value = "  +1.25%  "
$ws.Range("D20").Value = "'11.87"
$ws.Range("E20").Value = "  -0.31%  "
$ws.Range("E21").Value = "  +0.64%  "
$ws.Range("D22").Value = "'480.75"
$ws.Range("E22").Value = "  +1.90%  "
$ws.Range("D23").Value = "'4.95"
$ws.Range("E23").Value = "  -0.37%  "
$ws.Range("D24").Value = "'14.36"
$ws.Range("E24").Value = "  +4.61%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").Value = "'4.11"
$ws.Range("E25").Value = "  -0.31%  "
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").Value = "'89.51"
$ws.Range("E26").Value = "  +3.44%  "
$ws.Range("E27").Value = "  +2.01%  "
$ws.Range("D28").Value = "'10.65"
$ws.Range("E28").Value = "  -1.81%  "
$ws.Range("D30").Value = "'31.31"
$ws.Range("E30").Value = "  +2.21%  "
$ws.Range("D31").Value = "'6.57"
$ws.Range("E31").Value = "  -2.16%  "
$ws.Range("D32").Value = "'11.53"
$ws.Range("E32").Value = "  -0.15%  "
$ws.Range("E33").Value = "  +4.58%  "
$ws.Range("D34").Value = "'576.98"
$ws.Range("E34").Value = "  -0.55%  "
$ws.Range("E35").Value = "  -0.70%  "
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("D37").Value = "'3.65"
$ws.Range("E37").Value = "  +5.32%  "
$ws.Range("E38").Value = "  +0.76%  "
$ws.Range("D39").Value = "'35.89"
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("D40").Value = "'0.375"
$ws.Range("E40").Value = "  +0.71%  "
$ws.Range("D41").Value = "0.0₃0741"
$ws.Range("E41").Value = "  -1.77%  "
$ws.Range("D42").Value = "3.096.86"
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("D43").Value = "'2.81"
$ws.Range("E43").Value = "  -2.24%  "
$ws.Range("D44").Value = "'0.0418"
$ws.Range("E44").Value = "  +1.07%  "
$ws.Range("E45").Value = "  +0.33%  "
$ws.Range("D46").Value = "'3.17"
$ws.Range("E46").Value = "  -0.95%  "
$ws.Range("E47").Value = "  -3.35%  "
$ws.Range("D48").Value = "'0.998"
$ws.Range("E48").Value = "  -0.11%  "
$ws.Range("D49").Value = "'140.36"
$ws.Range("E49").Value = "  +2.24%  "
$ws.Range("E50").Value = "  +0.08%  "
$ws.Range("D51").Value = "'8.45"
$ws.Range("E51").Value = "  +0.84%  "

Write-Output "Applied cryptos update"
